$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $addr, $text) {
    # Force the cell to remain a text value even when $text parses as a
    # number (e.g. "0.9999"), matching the source workbook which stores
    # every Price/Volume cell as an inline string.
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" '26.049.35'
$ws.Range("E2").Value = '  +0.94%  '
Set-TextCell $ws "D3" '1.746.81'
Set-TextCell $ws "D4" '0.9999'
$ws.Range("E4").Value = '  -0.03%  '
Set-TextCell $ws "D5" '233.25'
$ws.Range("E5").Value = '  +2.42%  '
Set-TextCell $ws "D6" '1.000'
$ws.Range("E6").Value = '  +0.03%  '
Set-TextCell $ws "D7" '0.5274'
$ws.Range("E7").Value = '  +1.91%  '
Set-TextCell $ws "D8" '0.2766'
$ws.Range("E8").Value = '  +1.51%  '
Set-TextCell $ws "D9" '40.14'
$ws.Range("E9").Value = '  +4.27%  '
Set-TextCell $ws "D10" '0.06186'
$ws.Range("E10").Value = '  +1.55%  '
Set-TextCell $ws "D11" '1.753.95'
$ws.Range("E11").Value = '  +0.88%  '
Set-TextCell $ws "D12" '0.07205'
$ws.Range("E12").Value = '  +2.84%  '
Set-TextCell $ws "D13" '15.28'
$ws.Range("E13").Value = '  +0.32%  '
Set-TextCell $ws "D14" '0.6407'
$ws.Range("E14").Value = '  +1.66%  '
Set-TextCell $ws "D15" '4.590'
$ws.Range("E15").Value = '  +2.19%  '
Set-TextCell $ws "D16" '78.32'
$ws.Range("E16").Value = '  +2.64%  '
Set-TextCell $ws "D17" '1.000'
$ws.Range("E17").Value = '  +0.09%  '
Set-TextCell $ws "D18" '0.9998'
$ws.Range("E18").Value = '  -0.06%  '
Set-TextCell $ws "D19" '25.980.22'
$ws.Range("E19").Value = '  +0.62%  '
Set-TextCell $ws "D20" '11.58'
$ws.Range("E20").Value = '  +1.42%  '
Set-TextCell $ws "D21" '0.000006745'
$ws.Range("E21").Value = '  +2.34%  '
Set-TextCell $ws "D22" '1.977.42'
$ws.Range("E22").Value = '  +0.81%  '
Set-TextCell $ws "D23" '4.325'
$ws.Range("E23").Value = '  +6.88%  '
Set-TextCell $ws "D24" '8.817'
$ws.Range("E24").Value = '  +4.23%  '
Set-TextCell $ws "D25" '5.197'
$ws.Range("E25").Value = '  +1.90%  '
Set-TextCell $ws "D26" '139.46'
Set-TextCell $ws "D27" '1.521'
$ws.Range("E27").Value = '  +1.07%  '
Set-TextCell $ws "D28" '15.27'
$ws.Range("E28").Value = '  +1.91%  '
Set-TextCell $ws "D29" '1.806'
$ws.Range("E29").Value = '  -0.67%  '
Set-TextCell $ws "D30" '104.45'
$ws.Range("E30").Value = '  +1.90%  '
Set-TextCell $ws "D31" '0.08287'
$ws.Range("E31").Value = '  -0.44%  '
Set-TextCell $ws "D32" '3.767'
$ws.Range("E32").Value = '  +4.12%  '
Set-TextCell $ws "D33" '3.663'
$ws.Range("E33").Value = '  +8.77%  '
Set-TextCell $ws "D34" '0.04522'
$ws.Range("E34").Value = '  +2.43%  '
Set-TextCell $ws "D35" '2.636'
$ws.Range("E35").Value = '  +1.20%  '
Set-TextCell $ws "D36" '0.9989'
$ws.Range("E36").Value = '  +2.90%  '
Set-TextCell $ws "D37" '0.6309'
$ws.Range("E37").Value = '  +5.95%  '
Set-TextCell $ws "D38" '2.696'
$ws.Range("E38").Value = '  +0.51%  '
Set-TextCell $ws "D39" '0.01594'
$ws.Range("E39").Value = '  +2.29%  '
Set-TextCell $ws "D40" '1.928'
$ws.Range("E40").Value = '  -0.55%  '
Set-TextCell $ws "D41" '0.9996'
$ws.Range("E41").Value = '  +0.11%  '
Set-TextCell $ws "D42" '98.17'
$ws.Range("E42").Value = '  -3.49%  '
Set-TextCell $ws "D43" '0.3896'
$ws.Range("E43").Value = '  +2.58%  '
Set-TextCell $ws "D44" '0.7366'
$ws.Range("E44").Value = '  +1.65%  '
Set-TextCell $ws "D45" '5.039'
$ws.Range("E45").Value = '  +3.62%  '
Set-TextCell $ws "D46" '0.1141'
$ws.Range("E46").Value = '  +3.91%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws "D47" '6.308'
$ws.Range("E47").Value = '  +1.98%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws "D48" '0.05341'
$ws.Range("E48").Value = '  -2.58%  '
Set-TextCell $ws "D49" '53.98'
$ws.Range("E49").Value = '  +4.16%  '
Set-TextCell $ws "D50" '30.53'
$ws.Range("E50").Value = '  +2.46%  '
Set-TextCell $ws "D51" '7.636'
$ws.Range("E51").Value = '  +3.76%  '
